# edit.ps1 — applies the "feat: add 2022-Q3 data" change:
#  1. Insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q2"),
#     populated with the fund-holding breakdown for that quarter.
#  2. Insert a new top data row in "总计" for "2022-Q3" (49 holdings, 27.7
#     billion yuan), pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

function Set-HeaderLikeStyle {
    param($rng)
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $rng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $rng.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $rng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $rng.Borders.Item(11).LineStyle = 1  # xlInsideVertical
    $rng.Borders.Item(12).LineStyle = 1  # xlInsideHorizontal
}

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q3" sheet (inserted before the current #2 sheet,
#    i.e. right after "总计", matching the sheet order in the diff).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$q3.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
Set-HeaderLikeStyle -rng $q3.Range("B1:H1")

$q3Data = @(
    @("0","159941","广发纳斯达克100ETF（QDII）","106.15","91.14","4.52","4.7980","4"),
    @("1","513100","国泰纳斯达克100（QDII-ETF）","46.54","91.35","4.34","2.0198","4"),
    @("2","011421","广发全球科技三个月定期开放混合（QDII）美元 A","22.73","90.00","8.83","2.0071","1"),
    @("3","011420","广发全球科技三个月定期开放混合（QDII）人民币 A","22.73","90.00","8.83","2.0071","1"),
    @("4","000906","广发全球精选股票（QDII）美元现汇","21.88","79.27","8.89","1.9451","1"),
    @("5","270023","广发全球精选股票（QDII）","21.88","79.27","8.89","1.9451","1"),
    @("6","513500","博时标普500ETF（QDII）","71.37","96.44","2.26","1.6130","4"),
    @("7","040047","华安纳斯达克100指数（QDII）美元现钞A","24.52","92.09","4.55","1.1157","4"),
    @("8","040048","华安纳斯达克100指数（QDII）美元现汇A","24.52","92.09","4.55","1.1157","4"),
    @("9","040046","华安纳斯达克100指数（QDII）人民币A","22.21","92.09","4.55","1.0106","4"),
    @("10","002423","华宝标普美国品质消费股票（LOF）美元","3.59","94.37","20.52","0.7367","2"),
    @("11","160213","国泰纳斯达克100指数（QDII）","15.14","85.81","4.47","0.6768","4"),
    @("12","000041","华夏全球精选股票（QDII）","18.51","89.41","3.37","0.6238","6"),
    @("13","000834","大成纳斯达克100指数（QDII）","14.15","85.22","4.19","0.5929","4"),
    @("14","162415","华宝标普美国品质消费股票（LOF）人民币A","2.86","94.37","20.52","0.5869","2"),
    @("15","000043","嘉实美国成长股票（QDII）人民币","12.41","92.80","4.26","0.5287","5"),
    @("16","000044","嘉实美国成长股票（QDII）美元现汇","12.41","92.80","4.26","0.5287","5"),
    @("17","011423","广发全球科技三个月定期开放混合（QDII）美元 C","5.27","90.00","8.83","0.4653","1"),
    @("18","011422","广发全球科技三个月定期开放混合（QDII）人民币 C","5.27","90.00","8.83","0.4653","1"),
    @("19","012208","华夏港股前沿经济混合（QDII）A","9.67","89.48","4.29","0.4148","6"),
    @("20","001668","汇添富全球移动互联灵活配置混合（QDII）A","12.06","90.88","2.94","0.3546","6"),
    @("21","003722","易方达纳斯达克100指数美元（QDII-LOF）A","7.72","90.67","4.34","0.3350","4"),
    @("22","161130","易方达纳斯达克100指数人民币（QDII-LOF）","7.72","90.67","4.34","0.3350","4"),
    @("23","118002","易方达标普全球高端消费品指数增强A（QDII）人民币","1.85","93.04","9.14","0.1691","3"),
    @("24","000593","易方达标普全球高端消费品指数增强（QDII）美元现汇","1.85","93.04","9.14","0.1691","3"),
    @("25","005676","易方达标普全球高端消费品指数增强C（QDII）人民币","1.85","93.04","9.14","0.1691","3"),
    @("26","009975","华宝标普美国品质消费股票（LOF）人民币C","0.73","94.37","20.52","0.1498","2"),
    @("27","014978","华安纳斯达克100指数（QDII）人民币C","2.31","92.09","4.55","0.1051","4"),
    @("28","161125","易方达标普500指数（QDII-LOF）人民币","4.74","90.72","2.14","0.1014","4"),
    @("29","012860","易方达标普500指数（QDII-LOF）人民币 C","4.74","90.72","2.14","0.1014","4"),
    @("30","003718","易方达标普500指数（QDII-LOF）美元A","4.66","90.72","2.14","0.0997","4"),
    @("31","457001","国富亚洲机会股票（QDII）","3.80","83.80","2.32","0.0882","10"),
    @("32","159632","华安纳斯达克100ETF（QDII）","1.51","89.05","4.45","0.0672","4"),
    @("33","005698","华夏全球科技先锋混合（QDII）","0.59","86.79","8.98","0.0530","2"),
    @("34","161620","融通核心价值混合（QDII）A","0.55","57.96","9.53","0.0524","1"),
    @("35","013329","嘉实全球价值股票（QDII）美元现汇","1.68","90.63","1.71","0.0287","5"),
    @("36","013328","嘉实全球价值股票（QDII）人民币","1.68","90.63","1.71","0.0287","5"),
    @("37","006555","浦银安盛全球智能科技股票（QDII）A","0.25","84.65","7.40","0.0185","2"),
    @("38","015205","银华全球新能源车量化优选股票（QDII）C","0.24","90.37","7.01","0.0168","3"),
    @("39","012209","华夏港股前沿经济混合（QDII）C","0.35","89.48","4.29","0.0150","6"),
    @("40","015204","银华全球新能源车量化优选股票（QDII）A","0.19","90.37","7.01","0.0133","3"),
    @("41","159612","国泰标普500ETF（QDII）","0.55","91.40","1.98","0.0109","4"),
    @("42","012871","易方达纳斯达克100指数美元（QDII-LOF）C","0.18","90.67","4.34","0.0078","4"),
    @("43","012870","易方达纳斯达克100指数人民币（QDII-LOF）C","0.18","90.67","4.34","0.0078","4"),
    @("44","012861","易方达标普500指数（QDII-LOF）美元 C","0.08","90.72","2.14","0.0017","4"),
    @("45","015203","汇添富全球移动互联灵活配置混合（QDII）D","0.04","90.88","2.94","0.0012","6"),
    @("46","014127","融通核心价值混合（QDII）C","0.01","57.96","9.53","0.0010","1"),
    @("47","014002","浦银安盛全球智能科技股票（QDII）C","0.01","84.65","7.40","0.0007","2"),
    @("48","015202","汇添富全球移动互联灵活配置混合（QDII）C","0.01","90.88","2.94","0.0003","6")
)

foreach ($row in $q3Data) {
    $r = [int]$row[0] + 2
    $aCell = $q3.Cells.Item($r, 1)
    $aCell.Value = [int]$row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.Item(7).LineStyle = 1
    $aCell.Borders.Item(8).LineStyle = 1
    $aCell.Borders.Item(9).LineStyle = 1
    $aCell.Borders.Item(10).LineStyle = 1

    $bCell = $q3.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $q3.Cells.Item($r, 3).Value = $row[2]

    foreach ($colIdx in 4,5,6,7) {
        $cell = $q3.Cells.Item($r, $colIdx)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$colIdx - 1]
    }

    $q3.Cells.Item($r, 8).Value = [int]$row[7]
}

# ---------------------------------------------------------------------------
# 2. Insert the "2022-Q3" row at the top of the "总计" data and renumber.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

$summary.Range("B2:D2").Style = "Normal"
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 49
$summary.Cells.Item(2, 4).Value = 27.7

$aTop = $summary.Cells.Item(2, 1)
$aTop.Value = 0
$aTop.Font.Bold = $true
$aTop.HorizontalAlignment = -4108
$aTop.VerticalAlignment = -4160
$aTop.Borders.Item(7).LineStyle = 1
$aTop.Borders.Item(8).LineStyle = 1
$aTop.Borders.Item(9).LineStyle = 1
$aTop.Borders.Item(10).LineStyle = 1

# Renumber the remaining rows (0,1,2,3) to keep the sequential index in col A.
for ($r = 3; $r -le 5; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}
